# Apply scheduled market-data refresh to the per-job Leve profit tables.
# Each row's currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and LeveProfit(NQ/HQ)
# columns (H:N) are recomputed from the latest market board data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15: Morning Glass of Ether (Ether)
$ws.Range("H15").Value = 2175.4658
$ws.Range("I15").Value = 2175.4658
$ws.Range("K15").Value = 6526.3974
$ws.Range("M15").Value = -6357.3974
# Row 101: Edge of the Arcane (Cunning Craftsman's Tea)
$ws.Range("H101").Value = 1076
$ws.Range("I101").Value = 1003.4
$ws.Range("J101").Value = 1166.75
$ws.Range("K101").Value = 3010.2
$ws.Range("L101").Value = 3500.25
$ws.Range("M101").Value = -1388.2
$ws.Range("N101").Value = -6744.25
# Row 135: For Tired Minds (Grade 1 Gemsap of Intelligence)
$ws.Range("H135").Value = 1069.15
$ws.Range("I135").Value = 1010.1667
$ws.Range("J135").Value = 1600
$ws.Range("K135").Value = 9091.5003
$ws.Range("L135").Value = 14400
$ws.Range("M135").Value = -6556.5003
$ws.Range("N135").Value = -19470
# Row 137: Cutting Edge of Culinary Quality (Magnesia Whetstone)
$ws.Range("H137").Value = 2971053.5
$ws.Range("I137").Value = 276967.34
$ws.Range("K137").Value = 830902.02
$ws.Range("M137").Value = -828352.02
# Row 138: All-night Crafting (Cunning Craftsman's Tisane)
$ws.Range("H138").Value = 4389.94
$ws.Range("I138").Value = 2666.5
$ws.Range("K138").Value = 7999.5
$ws.Range("M138").Value = -2859.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust (Steel Ingot)
$ws.Range("H32").Value = 14780341
$ws.Range("I32").Value = 15257975
$ws.Range("K32").Value = 15257975
$ws.Range("M32").Value = -15257688
# Row 61: Dealing with the Tough Stuff (Cobalt Ingot)
$ws.Range("H61").Value = 5813.7144
$ws.Range("I61").Value = 5722.4614
$ws.Range("J61").Value = 7000
$ws.Range("K61").Value = 5722.4614
$ws.Range("L61").Value = 7000
$ws.Range("M61").Value = -5510.4614
$ws.Range("N61").Value = -7424
# Row 74: As the Bolt Flies (Titanium Nugget)
$ws.Range("H74").Value = 3416.3635
$ws.Range("I74").Value = 3358
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 3358
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -2484
$ws.Range("N74").Value = -5748
# Row 77: Heavy Metal Banned (L) (Titanium Nugget)
$ws.Range("H77").Value = 3416.3635
$ws.Range("I77").Value = 3358
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 16790
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -12422
$ws.Range("N77").Value = -28736
# Row 102: Smells of Rich Tama-hagane (Tama-hagane Ingot)
$ws.Range("H102").Value = 1886.1538
$ws.Range("I102").Value = 1352
$ws.Range("J102").Value = 3666.6667
$ws.Range("K102").Value = 1352
$ws.Range("L102").Value = 3666.6667
$ws.Range("M102").Value = 270
$ws.Range("N102").Value = -6910.6667
# Row 132: Don't Bore Me, Ore Me (Mountain Chromite Ingot)
$ws.Range("H132").Value = 5268.909
$ws.Range("I132").Value = 5117.6665
$ws.Range("J132").Value = 5949.5
$ws.Range("K132").Value = 15352.9995
$ws.Range("L132").Value = 17848.5
$ws.Range("M132").Value = -12822.9995
$ws.Range("N132").Value = -22908.5
# Row 136: Metal with Mettle (Cobalt Tungsten Ingot)
$ws.Range("H136").Value = 5813.7144
$ws.Range("I136").Value = 5722.4614
$ws.Range("J136").Value = 7000
$ws.Range("K136").Value = 17167.3842
$ws.Range("L136").Value = 21000
$ws.Range("M136").Value = -14617.3842
$ws.Range("N136").Value = -26100
# Row 138: Don't Ask about the Rivets (Titanium Gold Helm of Casting)
$ws.Range("H138").Value = 99979
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin (Adamantite Nugget)
$ws.Range("H86").Value = 3249.125
$ws.Range("I86").Value = 2806.6924
$ws.Range("J86").Value = 5166.3335
$ws.Range("K86").Value = 2806.6924
$ws.Range("L86").Value = 5166.3335
$ws.Range("M86").Value = -1683.6924
$ws.Range("N86").Value = -7412.3335
# Row 89: Piercing Eyes Deserve Piercing Shafts (L) (Adamantite Nugget)
$ws.Range("H89").Value = 3249.125
$ws.Range("I89").Value = 2806.6924
$ws.Range("J89").Value = 5166.3335
$ws.Range("K89").Value = 14033.462
$ws.Range("L89").Value = 25831.6675
$ws.Range("M89").Value = -8417.462
$ws.Range("N89").Value = -37063.6675
# Row 105: Ingot to Wing It (Molybdenum Ingot)
$ws.Range("H105").Value = 3568.9167
$ws.Range("I105").Value = 3591.125
$ws.Range("K105").Value = 3591.125
$ws.Range("M105").Value = -1844.125

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found (Walnut Lumber)
$ws.Range("H31").Value = 3461.4707
$ws.Range("I31").Value = 1829.1936
$ws.Range("J31").Value = 5991.5
$ws.Range("K31").Value = 1829.1936
$ws.Range("L31").Value = 5991.5
$ws.Range("M31").Value = -1534.1936
$ws.Range("N31").Value = -6581.5
# Row 34: Armoires of the Rich and Famous (Walnut Lumber)
$ws.Range("H34").Value = 3461.4707
$ws.Range("I34").Value = 1829.1936
$ws.Range("J34").Value = 5991.5
$ws.Range("K34").Value = 1829.1936
$ws.Range("L34").Value = 5991.5
$ws.Range("M34").Value = -1627.1936
$ws.Range("N34").Value = -6395.5
# Row 132: Hull Lotta Damage (Ginseng Lumber)
$ws.Range("H132").Value = 1551.45
$ws.Range("I132").Value = 1165.1875
$ws.Range("J132").Value = 3096.5
$ws.Range("K132").Value = 3495.5625
$ws.Range("L132").Value = 9289.5
$ws.Range("M132").Value = -965.5625
$ws.Range("N132").Value = -14349.5
# Row 134: Wood You Be Quiet (Ceiba Lumber)
$ws.Range("H134").Value = 2053.6
$ws.Range("I134").Value = 2014.56
$ws.Range("J134").Value = 2248.8
$ws.Range("K134").Value = 6043.68
$ws.Range("L134").Value = 6746.400000000001
$ws.Range("M134").Value = -3508.68
$ws.Range("N134").Value = -11816.4

$ws = $wb.Worksheets.Item("CUL")
# Row 68: Such a Butter Face (Fermented Butter)
$ws.Range("H68").Value = 3513
$ws.Range("I68").Value = 2983.3333
$ws.Range("K68").Value = 8949.999899999999
$ws.Range("M68").Value = -8138.999899999999
# Row 71: No Margarine of Error (L) (Fermented Butter)
$ws.Range("H71").Value = 3513
$ws.Range("I71").Value = 2983.3333
$ws.Range("K71").Value = 26849.9997
$ws.Range("M71").Value = -22793.9997

$ws = $wb.Worksheets.Item("LTW")
# Row 127: Loyal Turncoat (Saigaskin Coat of Fending)
$ws.Range("H127").Value = 23181.773
$ws.Range("J127").Value = 23181.773
$ws.Range("L127").Value = 23181.773
$ws.Range("N127").Value = -33101.773
# Row 132: Tenets of Tanning (Silver Lobo Leather)
$ws.Range("H132").Value = 5934.85
$ws.Range("I132").Value = 5733.923
$ws.Range("J132").Value = 6308
$ws.Range("K132").Value = 17201.769
$ws.Range("L132").Value = 18924
$ws.Range("M132").Value = -14671.769
$ws.Range("N132").Value = -23984

$ws = $wb.Worksheets.Item("WVR")
# Row 46: Crunching the Numbers (Linen Hat)
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
# Row 100: Of Great Import (Kudzu Thread)
$ws.Range("H100").Value = 453.25
$ws.Range("I100").Value = 223.25
$ws.Range("K100").Value = 446.5
$ws.Range("M100").Value = 94.5
# Row 122: Heavy Armoire (Dark Hempen Cloth)
$ws.Range("H122").Value = 4285.606
$ws.Range("I122").Value = 4201.1113
$ws.Range("J122").Value = 4665.8335
$ws.Range("K122").Value = 12603.3339
$ws.Range("L122").Value = 13997.5005
$ws.Range("M122").Value = -10153.3339
$ws.Range("N122").Value = -18897.5005
# Row 132: Comfy Cabins (Snow Cotton Cloth)
$ws.Range("H132").Value = 3377.6
$ws.Range("I132").Value = 3066.7693
$ws.Range("J132").Value = 3954.8572
$ws.Range("K132").Value = 9200.3079
$ws.Range("L132").Value = 11864.5716
$ws.Range("M132").Value = -6670.3079
$ws.Range("N132").Value = -16924.5716
# Row 134: Cloth for Canvas (Mountain Linen)
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
# Row 136: Weaving the Envelope (Sarcenet Cloth)
$ws.Range("H136").Value = 87998.836
$ws.Range("I136").Value = 3766.5
$ws.Range("J136").Value = 172231.17
$ws.Range("K136").Value = 11299.5
$ws.Range("L136").Value = 516693.51
$ws.Range("M136").Value = -8749.5
$ws.Range("N136").Value = -521793.51
